$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Normalize shorthand shift-code labels that no longer carry their
# qualifier suffix in the updated schedule.
[void]$ws.Cells.Replace("10 - 3 (M)", "10 - 3")
[void]$ws.Cells.Replace("9 - 2 (M)", "9 - 2")
[void]$ws.Cells.Replace("6 - 1 (S)", "6 - 1")

# Remove the old legend/key block at the bottom of the sheet (rows 62-66).
$ws.Rows("62:66").Delete()

# Leave the selection where the author last left it while reviewing the sheet.
[void]$ws.Range("D70").Select()
